$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "type.php" Type-model block (and everything below it, down through the
# cart-relation functions block) lives in columns H:M and needs to shift
# down by 2 rows, starting at row 11, to make room for two new "Functions :"
# entries (convertDH / convertHD) that belong to the style.php block above it.
#
# Columns A:F and O:T are untouched - only H:M moves. We walk bottom-to-top
# so we never clobber a source row before it has been read.

for ($r = 29; $r -ge 11; $r--) {
    $srcRow = $r
    $dstRow = $r + 2
    $srcRange = $ws.Range("H$srcRow" + ":M$srcRow")
    $dstRange = $ws.Range("H$dstRow" + ":M$dstRow")

    $hasContent = $false
    for ($c = 8; $c -le 13; $c++) {
        $cell = $ws.Cells.Item($srcRow, $c)
        if ($cell.Value -ne $null) { $hasContent = $true }
    }

    if ($hasContent) {
        $srcRange.Copy($dstRange)
    } else {
        $dstRange.Clear()
    }
}

# Row 11 becomes a blank spacer row styled like the new rows above it (s=4)
# instead of its old contents (it used to hold the "type.php" header, s=2).
$ws.Range("H11:M11").ClearContents()
$ws.Range("H9:M9").Copy()
$ws.Range("H11:M11").PasteSpecial(-4122)

# Row 12 had the "Objet : / Type" header - that row is now fully empty in
# H:M (no cells at all).
$ws.Range("H12:M12").Clear()

# Give row 10 (brand new) the same fill style as row 9 (s=4) before filling
# it in - row 10 previously had no H:M cells whatsoever.
$ws.Range("H9:M9").Copy()
$ws.Range("H10:M10").PasteSpecial(-4122)

# Populate the two new Functions rows.
$ws.Range("I9").Value = "convertDH"
$ws.Range("J9").Value = "decimal"
$ws.Range("K9").Value = "Converti le décimal en hexa"

$ws.Range("I10").Value = "convertHD"
$ws.Range("J10").Value = "hexa"
$ws.Range("K10").Value = "Converti le hexa en decimal"

$excel.CutCopyMode = $false

# Match the author's view state from the saved file (scrolled right a bit,
# selection resting on K11).
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("K11").Select()
